$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows for courses that were removed:
#  - Row 24: PROYECTO DE SISTEMAS ROBUSTOS, PARALELOS Y DISTRIBUIDOS
#  - Row 23: PROYECTO DE GESTION DE LA TECNOLOGIA DE INFORMACION
#  - Row 9:  COMPUTO FLEXIBLE (SOFTCOMPUTING)
# Deleted from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(9).Delete()
